$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the hyperlink that was attached to A52 (the eurostat database URL).
# The text itself moves to a plain (non-hyperlink) cell further down.
foreach ($h in $ws.Hyperlinks) {
    $h.Delete()
}

# --- Row 51: used to hold "SBS Main Indicators..." -> now fully blank (no value, no format) ---
$ws.Range("A51").Clear()

# --- Row 52: used to be the hyperlinked URL -> now the "SBS Main Indicators..." text,
# formatted like the surrounding "source" text (italic, default color, no underline) ---
$c52 = $ws.Range("A52")
$c52.Value = "SBS Main Indicators, Annual enterprise statistics by size class for special aggregates of activities (NACE Rev. 2)"
$c52.Font.Underline = $false
$c52.Font.ColorIndex = 0
$c52.Font.Italic = $true

# --- Row 53: already an empty string with the "source" style -> leave untouched ---

# --- Row 54 (new row): the eurostat database URL as plain text, "source" style (italic) ---
$c54 = $ws.Range("A54")
$c54.Value = "http://epp.eurostat.ec.europa.eu/portal/page/portal/european_business/data/database"
$c54.Font.Italic = $true

# --- Row 56: used to hold "SBS Eurostat" -> now fully blank (no value, no format) ---
$ws.Range("A56").Clear()

# --- Row 57: now holds "SBS Eurostat" with the bold "title" style ---
$c57 = $ws.Range("A57")
$c57.Value = "SBS Eurostat"
$c57.Font.Bold = $true
$c57.Font.Italic = $false
$c57.Font.Underline = $false

# --- Row 58 (new row): also "SBS Eurostat", but with the italic "source" style ---
$c58 = $ws.Range("A58")
$c58.Value = "SBS Eurostat"
$c58.Font.Italic = $true
$c58.Font.Bold = $false
